# feat: default format for datasets now jsonl
#
# Renames the example dataset's header columns from the old
# query/positive/negative triad to the jsonl-style prompt/chosen/rejected
# naming, and gives the table body (everything below/right of the header
# row) an explicit white interior fill to match the refreshed template
# styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (shared strings: query -> prompt, positive -> chosen,
# negative -> rejected).
$ws.Range("A1").Value = "prompt"
$ws.Range("B1").Value = "chosen"
$ws.Range("C1").Value = "rejected"

# Apply an explicit white fill across the table body cells (the header row
# A1:C1 already carries its own accent fill and is left untouched).
# 16777215 = 0xFFFFFF -> RGB(255, 255, 255) in the BGR-packed OLE COLORREF
# format Excel's Interior.Color expects.
$white = 16777215
$ws.Range("D1:E1").Interior.Color = $white
$ws.Range("A2:E10").Interior.Color = $white
